$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Add the new header "l1" in column R, row 1
$ws.Range("R1").Value = "l1"

# Copy formatting (bold + border + alignment) from the existing Q1 header cell
$ws.Range("Q1").Copy() | Out-Null
$ws.Range("R1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill column R with 0 for every data row, matching columns P and Q
$dataRange = $ws.Range("R2:R" + $lastRow)
$dataRange.Value = 0
